$d = $word.ActiveDocument

function New-WordOpenXmlPackage($bodyInnerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step 1: drop the bookmark that wraps the first (title) paragraph ---
# The bookmarkStart/bookmarkEnd live outside any w:p's run content, so no
# Range-based rewrite of the paragraph touches them. Deleting the whole
# paragraph first collapses them down to zero-width stubs at offset 0,
# which a couple of zero-length Range deletes then consume.
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$titleRange.Delete()

$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# --- Step 2: re-create the title paragraph, styled "Title", split word by word ---
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()

$titleBody = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Fall</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Appeal</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">-</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">October</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">1976</w:t></w:r>' +
    '</w:p>'

$newTitlePara = $d.Paragraphs.Item(1)
$newTitlePara.Range.InsertXML((New-WordOpenXmlPackage $titleBody))

# --- Step 3: rewrite "By Dorothy Day" paragraph as the "Authors" block ---
$authorsBody = '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' +
    '</w:p>'

$authorsPara = $d.Paragraphs.Item(2)
$authorsPara.Range.InsertXML((New-WordOpenXmlPackage $authorsBody))
